# Add a new "Georgia" row to the county-level election results table,
# placed alphabetically between "Florida" (row 9) and "Iowa" (row 10).
# Every row from Iowa onward shifts down by one to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a blank worksheet row at row 10; this pushes Iowa..Wisconsin down
# by one row (rows 10-29 become rows 11-30).
$ws.Rows.Item(10).Insert()

# Grow the table (and its autofilter) so it covers the extra row.
$lo.Resize($ws.Range("A1:C30"))

# Populate the new Georgia row.
$ws.Range("A10").Value = "Georgia"
$ws.Range("B10").Value = "https://en.wikipedia.org/wiki/United_States_presidential_election_in_Georgia,_2016#Results_by_County"
$ws.Range("C10").Formula = "=""{ state: ""&CHAR(34)&Table1[[#This Row],[State]]&CHAR(34)&"",URL: ""&CHAR(34)&Table1[[#This Row],[URL]]&CHAR(34)&""},"""

# The last row was pushed past the table's original bounds during the
# Insert, so its formula was left referencing a calculated column
# outside the (not-yet-resized) table and not recalculated. Re-apply it
# now that the table covers A1:C30 so it resolves to the right value.
$ws.Range("C30").Formula = "=""{ state: ""&CHAR(34)&Table1[[#This Row],[State]]&CHAR(34)&"",URL: ""&CHAR(34)&Table1[[#This Row],[URL]]&CHAR(34)&""},"""

# Match the saved cursor position from the edit.
[void]$ws.Range("C10").Select()
